$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster = ECs)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.06449866666666666
$ws.Range("N2").Value = 0.193496
$ws.Range("O2").Value = 0.004525829983623641
$ws.Range("P2").Value = 0.004525829983623642
$ws.Range("Q2").Value = 0.004596647976888889
$ws.Range("R2").Value = 0.041369831792
$ws.Range("S2").Value = 0.004525829983623641
$ws.Range("T2").Value = 0.004525829983623642

# Row 3 (Target cluster = FAPs)
$ws.Range("O3").Value = 0.745188142173877
$ws.Range("P3").Value = 0.7451881421738772
$ws.Range("S3").Value = 0.745188142173877
$ws.Range("T3").Value = 0.7451881421738772

# Row 4 (Target cluster = MuSCs)
$ws.Range("M4").Value = 3.566885000000001
$ws.Range("N4").Value = 10.700655
$ws.Range("O4").Value = 0.2502860278424993
$ws.Range("P4").Value = 0.2502860278424993
$ws.Range("Q4").Value = 0.2542023822566667
$ws.Range("R4").Value = 2.28782144031
$ws.Range("S4").Value = 0.2502860278424993
$ws.Range("T4").Value = 0.2502860278424993
